$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 2165100
$ws.Range("E8").Value = 2078500
$ws.Range("F8").Value = 2259100
$ws.Range("G8").Value = 1972100
$ws.Range("H8").Value = 1967100
$ws.Range("I8").Value = 2043000
$ws.Range("J8").Value = 2402500

$ws.Range("D9").Value = 1285700
$ws.Range("E9").Value = 1276400
$ws.Range("F9").Value = 1468500
$ws.Range("G9").Value = 1326200
$ws.Range("H9").Value = 2727700
$ws.Range("I9").Value = 1381400
$ws.Range("J9").Value = 1576700

$ws.Range("D10").Value = 879300
$ws.Range("E10").Value = 802100
$ws.Range("F10").Value = 790600
$ws.Range("G10").Value = 645900
$ws.Range("H10").Value = -760600
$ws.Range("I10").Value = 661600
$ws.Range("J10").Value = 825700

$ws.Range("H12").Value = 23700
$ws.Range("I12").Value = 21800

$ws.Range("E14").Value = 4900
$ws.Range("F14").Value = 37700
$ws.Range("G14").Value = 48500
$ws.Range("H14").Value = 110900
$ws.Range("J14").Value = 3100

$ws.Range("D17").Value = 1756600
$ws.Range("E17").Value = 1749800
$ws.Range("F17").Value = 2036000
$ws.Range("G17").Value = 1833800
$ws.Range("H17").Value = 1896300
$ws.Range("I17").Value = 1845200
$ws.Range("J17").Value = 2032300

$ws.Range("D18").Value = 408400
$ws.Range("E18").Value = 328700
$ws.Range("F18").Value = 223100
$ws.Range("G18").Value = 138400
$ws.Range("H18").Value = 70700
$ws.Range("I18").Value = 197800
$ws.Range("J18").Value = 370200

$ws.Range("D20").Value = 3200
$ws.Range("G20").Value = 24300
$ws.Range("H20").Value = 24800
$ws.Range("I20").Value = 12100

$ws.Range("D21").Value = 524700
$ws.Range("E21").Value = 477000
$ws.Range("F21").Value = 486700
$ws.Range("G21").Value = 349300
$ws.Range("H21").Value = 287600
$ws.Range("I21").Value = 305000
$ws.Range("J21").Value = 463400

$ws.Range("D22").Value = 7400
$ws.Range("E22").Value = 8500
$ws.Range("F22").Value = 8600
$ws.Range("G22").Value = 9300
$ws.Range("H22").Value = 10700
$ws.Range("I22").Value = 11800
$ws.Range("J22").Value = 12900

$ws.Range("D23").Value = 404200
$ws.Range("E23").Value = 321100
$ws.Range("F23").Value = 214900
$ws.Range("G23").Value = 153300
$ws.Range("H23").Value = 84800
$ws.Range("I23").Value = 198100
$ws.Range("J23").Value = 361800

$ws.Range("D24").Value = 128400
$ws.Range("E24").Value = 86300
$ws.Range("F24").Value = 119700
$ws.Range("G24").Value = 63200
$ws.Range("H24").Value = 43600
$ws.Range("I24").Value = 76600
$ws.Range("J24").Value = 153100

$ws.Range("D26").Value = 275800
$ws.Range("E26").Value = 234800
$ws.Range("F26").Value = 95200
$ws.Range("G26").Value = 90100
$ws.Range("H26").Value = 41100
$ws.Range("I26").Value = 121500
$ws.Range("J26").Value = 208700

$ws.Range("D27").Value = 275800
$ws.Range("E27").Value = 234600
$ws.Range("F27").Value = 95100
$ws.Range("G27").Value = 89700
$ws.Range("H27").Value = 40600
$ws.Range("I27").Value = 119100
$ws.Range("J27").Value = 208000

$ws.Range("D32").Value = -3200
$ws.Range("G32").Value = -24300
$ws.Range("H32").Value = -24800
$ws.Range("I32").Value = -12100

$ws.Range("D33").Value = 275800
$ws.Range("E33").Value = 234600
$ws.Range("F33").Value = 95100
$ws.Range("G33").Value = 89700
$ws.Range("H33").Value = 40600
$ws.Range("I33").Value = 119100
$ws.Range("J33").Value = 208000

$ws.Range("D35").Value = 275800
$ws.Range("E35").Value = 234600
$ws.Range("F35").Value = 95100
$ws.Range("G35").Value = 89700
$ws.Range("H35").Value = 40600
$ws.Range("I35").Value = 119100
$ws.Range("J35").Value = 208000

$ws.Range("D41").Value = 1365700
$ws.Range("E41").Value = 1200900
$ws.Range("F41").Value = 1016900
$ws.Range("G41").Value = 532500
$ws.Range("H41").Value = 904400
$ws.Range("I41").Value = 575600
$ws.Range("J41").Value = 691100

$ws.Range("D43").Value = 241900
$ws.Range("E43").Value = 242200
$ws.Range("F43").Value = 228600
$ws.Range("G43").Value = 297100
$ws.Range("H43").Value = 560900
$ws.Range("I43").Value = 298900
$ws.Range("J43").Value = 304200

$ws.Range("D44").Value = 61800
$ws.Range("E44").Value = 67200
$ws.Range("F44").Value = 82900
$ws.Range("G44").Value = 116100
$ws.Range("H44").Value = 381900
$ws.Range("I44").Value = 238200
$ws.Range("J44").Value = 200000

$ws.Range("D45").Value = 99400
$ws.Range("E45").Value = 71300
$ws.Range("F45").Value = 65900
$ws.Range("G45").Value = 106200
$ws.Range("H45").Value = 345200
$ws.Range("I45").Value = 274800
$ws.Range("J45").Value = 268900

$ws.Range("D46").Value = 1768900
$ws.Range("E46").Value = 1581500
$ws.Range("F46").Value = 1394300
$ws.Range("G46").Value = 1052000
$ws.Range("H46").Value = 929900
$ws.Range("I46").Value = 1387500
$ws.Range("J46").Value = 1464200

$ws.Range("D47").Value = 231500
$ws.Range("E47").Value = 227400
$ws.Range("F47").Value = 241400
$ws.Range("G47").Value = 240700
$ws.Range("H47").Value = 264400
$ws.Range("I47").Value = 24300
$ws.Range("J47").Value = 23600

$ws.Range("D48").Value = 714900
$ws.Range("E48").Value = 683400
$ws.Range("F48").Value = 725600
$ws.Range("G48").Value = 716500
$ws.Range("H48").Value = 1424000
$ws.Range("I48").Value = 566400
$ws.Range("J48").Value = 562700

$ws.Range("D49").Value = 333300
$ws.Range("E49").Value = 314400
$ws.Range("F49").Value = 356800
$ws.Range("G49").Value = 551800
$ws.Range("H49").Value = 1095300
$ws.Range("I49").Value = 580000
$ws.Range("J49").Value = 570900

$ws.Range("D52").Value = 233900
$ws.Range("E52").Value = 241000
$ws.Range("F52").Value = 248700
$ws.Range("G52").Value = 255800
$ws.Range("H52").Value = 636300
$ws.Range("I52").Value = 361300
$ws.Range("J52").Value = 343700

$ws.Range("D54").Value = 3282500
$ws.Range("E54").Value = 3047800
$ws.Range("F54").Value = 2966800
$ws.Range("G54").Value = 2816800
$ws.Range("H54").Value = 2717400
$ws.Range("I54").Value = 2919400
$ws.Range("J54").Value = 2965200

$ws.Range("D57").Value = 96100
$ws.Range("E57").Value = 75300
$ws.Range("F57").Value = 75800
$ws.Range("G57").Value = 99100
$ws.Range("H57").Value = 188900
$ws.Range("I57").Value = 130600
$ws.Range("J57").Value = 147300

$ws.Range("D58").Value = 107600
$ws.Range("E58").Value = 95900
$ws.Range("F58").Value = 81500
$ws.Range("G58").Value = 54300
$ws.Range("H58").Value = 135500
$ws.Range("I58").Value = 107100
$ws.Range("J58").Value = 88200

$ws.Range("D59").Value = 422700
$ws.Range("E59").Value = 342700
$ws.Range("F59").Value = 398600
$ws.Range("G59").Value = 313000
$ws.Range("H59").Value = 521500
$ws.Range("I59").Value = 300300
$ws.Range("J59").Value = 378300

$ws.Range("D60").Value = 626400
$ws.Range("E60").Value = 513900
$ws.Range("F60").Value = 555900
$ws.Range("G60").Value = 466500
$ws.Range("H60").Value = 436200
$ws.Range("I60").Value = 538000
$ws.Range("J60").Value = 613700

$ws.Range("D61").Value = 133300
$ws.Range("E61").Value = 310200
$ws.Range("F61").Value = 371200
$ws.Range("G61").Value = 301900
$ws.Range("H61").Value = 320100
$ws.Range("I61").Value = 204200
$ws.Range("J61").Value = 269400

$ws.Range("D62").Value = 221800
$ws.Range("E62").Value = 97600
$ws.Range("F62").Value = 109900
$ws.Range("G62").Value = 73300
$ws.Range("H62").Value = 196500
$ws.Range("I62").Value = 134200
$ws.Range("J62").Value = 131900

$ws.Range("D66").Value = 988300
$ws.Range("E66").Value = 928400
$ws.Range("F66").Value = 1043600
$ws.Range("G66").Value = 848000
$ws.Range("H66").Value = 835400
$ws.Range("I66").Value = 881600
$ws.Range("J66").Value = 1017400

$ws.Range("D72").Value = 1380100
$ws.Range("E72").Value = 1191100
$ws.Range("F72").Value = 992600
$ws.Range("G72").Value = 926400
$ws.Range("H72").Value = 1882100
$ws.Range("I72").Value = 1031400
$ws.Range("J72").Value = 975000

$ws.Range("D76").Value = 2294200
$ws.Range("E76").Value = 2119300
$ws.Range("F76").Value = 1923300
$ws.Range("G76").Value = 1968800
$ws.Range("H76").Value = 1881900
$ws.Range("I76").Value = 2037800
$ws.Range("J76").Value = 1947700

$ws.Range("D81").Value = 275800
$ws.Range("E81").Value = 234600
$ws.Range("F81").Value = 95100
$ws.Range("G81").Value = 89700
$ws.Range("H81").Value = 40600
$ws.Range("I81").Value = 119100
$ws.Range("J81").Value = 208000

$ws.Range("D83").Value = 112900
$ws.Range("E83").Value = 147200
$ws.Range("F83").Value = 262900
$ws.Range("G83").Value = 186500
$ws.Range("H83").Value = 191900
$ws.Range("I83").Value = 95000
$ws.Range("J83").Value = 88600

$ws.Range("D89").Value = 488000
$ws.Range("E89").Value = 395600
$ws.Range("F89").Value = 644900
$ws.Range("G89").Value = 409100
$ws.Range("H89").Value = 268600
$ws.Range("I89").Value = 92500
$ws.Range("J89").Value = 342800

$ws.Range("D91").Value = -159400
$ws.Range("E91").Value = -135300
$ws.Range("F91").Value = -172500
$ws.Range("G91").Value = -233000
$ws.Range("H91").Value = -667100
$ws.Range("I91").Value = -89700
$ws.Range("J91").Value = -83700

$ws.Range("D94").Value = -167100
$ws.Range("E94").Value = -123100
$ws.Range("F94").Value = -169500
$ws.Range("G94").Value = -221400
$ws.Range("H94").Value = -428600
$ws.Range("I94").Value = -104600
$ws.Range("J94").Value = -69100

$ws.Range("D96").Value = -86700
$ws.Range("E96").Value = -36000
$ws.Range("F96").Value = -28800
$ws.Range("G96").Value = -31900
$ws.Range("H96").Value = -52600
$ws.Range("I96").Value = -62500
$ws.Range("J96").Value = -51400

$ws.Range("D100").Value = -133800
$ws.Range("E100").Value = -85200
$ws.Range("F100").Value = -17000
$ws.Range("G100").Value = -61500
$ws.Range("H100").Value = 31200
$ws.Range("I100").Value = -111900
$ws.Range("J100").Value = -119800

$ws.Range("D101").Value = -8600
$ws.Range("F101").Value = -13200
$ws.Range("G101").Value = 6100
$ws.Range("H101").Value = 5600
$ws.Range("I101").Value = 8400

$ws.Range("D102").Value = 178500
$ws.Range("E102").Value = 188400
$ws.Range("F102").Value = 445200
$ws.Range("G102").Value = 132300
$ws.Range("H102").Value = -123400
$ws.Range("I102").Value = -115500
$ws.Range("J102").Value = 152900
